$d = $word.ActiveDocument

# Accept all tracked changes everywhere in the document (body + headers/footers).
$d.Revisions.AcceptAll()

# Temporarily stop tracking so the structural edits below (paragraph removal,
# bookmark relocation, date fix) are applied as plain content changes rather
# than becoming new tracked insertions/deletions.
$originalTrackRevisions = $d.TrackRevisions
$d.TrackRevisions = $false

# Remove the paragraph recommending an editor ("We suggest Tim F. Cooper as an
# editor for the paper, because of his previous work on complex adaptation and
# evolvability.") entirely.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("We suggest")) {
        $p.Range.Delete()
        $found = $true
        break
    }
}
Write-Host "Removed editor-recommendation paragraph: " $found

# The removed paragraph used to precede the "We hope you will find..."
# paragraph, which itself used to contain the (now accepted, mid-paragraph)
# "_GoBack" bookmark. Re-anchor that bookmark at the very start of the
# "We hope..." paragraph.
$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("We hope you will find")) {
        $r = $p.Range.Duplicate
        $r.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $r)
        $found2 = $true
        break
    }
}
Write-Host "Relocated _GoBack bookmark: " $found2

# Bump the header date from 27/4/2014 to 28/4/2014 (change the "7" to "8").
$hdr = $d.Sections(1).Headers(1)
$hrng = $hdr.Range.Duplicate
$hrng.Find.Execute("27/4/2014", $false, $false, $false, $false, $false, $true, 1, $false, "28/4/2014", 2)

# Restore the document's original track-changes setting.
$d.TrackRevisions = $originalTrackRevisions
